$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix date modified for serial no. 3 and 4 (rows 4 and 5) -> 2018-02-13
$ws.Range("B4").Value = 43144
$ws.Range("B5").Value = 43144

# New log entry for M1D4 (row 7)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 43144
$ws.Range("C7").Value = 0.82638888888888884
$ws.Range("D7").Value = "Team_04_M1_D4_Tool_Policy"
$ws.Range("E7").Value = "2100-Prerana"
$ws.Range("G7").Value = "Initial Version"

$ws.Range("G7").Select()
